# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 25 de Octubre de 2020 a las 09:00"

# --- Row 7: Rusia ---
$ws.Cells.Item(7,2).Value = 1513877
$ws.Cells.Item(7,3).Value = 16710
$ws.Cells.Item(7,4).Value = 1138522
$ws.Cells.Item(7,5).Value = 349305
$ws.Cells.Item(7,7).Value = 229
$ws.Cells.Item(7,8).Value = 26050

# --- Row 26: Ucrania ---
$ws.Cells.Item(26,2).Value = 343498
$ws.Cells.Item(26,3).Value = 6088
$ws.Cells.Item(26,4).Value = 141508
$ws.Cells.Item(26,5).Value = 195599
$ws.Cells.Item(26,7).Value = 102
$ws.Cells.Item(26,8).Value = 6391

# --- Row 60: Armenia ---
$ws.Cells.Item(60,2).Value = 77837
$ws.Cells.Item(60,3).Value = 2314
$ws.Cells.Item(60,4).Value = 50908
$ws.Cells.Item(60,5).Value = 25749
$ws.Cells.Item(60,7).Value = 23
$ws.Cells.Item(60,8).Value = 1180

# --- Row 66: Singapur ---
$ws.Cells.Item(66,2).Value = 57970
$ws.Cells.Item(66,3).Value = 5
$ws.Cells.Item(66,5).Value = 98

# --- Rows 88/89: Georgia inserted ahead of Australia (sst swap), values refreshed ---
$ws.Cells.Item(88,1).Value = "Georgia"
$ws.Cells.Item(88,2).Value = 28431
$ws.Cells.Item(88,3).Value = 1928
$ws.Cells.Item(88,4).Value = 10767
$ws.Cells.Item(88,5).Value = 17463
$ws.Cells.Item(88,7).Value = 8
$ws.Cells.Item(88,8).Value = 201

$ws.Cells.Item(89,1).Value = "Australia"
$ws.Cells.Item(89,2).Value = 27513
$ws.Cells.Item(89,3).Value = 14
$ws.Cells.Item(89,4).Value = 25181
$ws.Cells.Item(89,5).Value = 1427
$ws.Cells.Item(89,7).Value = 0
$ws.Cells.Item(89,8).Value = 905

# --- Rows 216/217: Islas Malvinas inserted ahead of Montserrat (sst swap), values refreshed ---
$ws.Cells.Item(216,1).Value = "Islas Malvinas"
$ws.Cells.Item(216,2).Value = 13
$ws.Cells.Item(216,3).Value = 0
$ws.Cells.Item(216,4).Value = 13
$ws.Cells.Item(216,5).Value = 0
$ws.Cells.Item(216,7).Value = 0
$ws.Cells.Item(216,8).Value = 0

$ws.Cells.Item(217,1).Value = "Montserrat"
$ws.Cells.Item(217,2).Value = 13
$ws.Cells.Item(217,3).Value = 0
$ws.Cells.Item(217,4).Value = 12
$ws.Cells.Item(217,5).Value = 0
$ws.Cells.Item(217,7).Value = 0
$ws.Cells.Item(217,8).Value = 1
